$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 43 ("「小さい」" entry) — causes all subsequent rows to shift up by one.
$ws.Rows(43).Delete()
